$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (symbol list refresh).
# Each entry is CellRef = NewTextValue. Values are written as text
# (matching the workbook's original inlineStr string storage for these
# columns), so NumberFormat is forced to "@" before assignment to stop
# Excel from auto-coercing numeric-/percent-looking strings into numbers.
$updates = [ordered]@{
    "D2"  = "273.90"
    "E2"  = "-1.75%"
    "D3"  = "26.72"
    "E3"  = "-1.93%"
    "D4"  = "4.863"
    "E4"  = "0.86%"
    "D5"  = "0.06322"
    "E5"  = "0.93%"
    "E6"  = "0.36%"
    "D7"  = "3.322"
    "E7"  = "1.74%"
    "D8"  = "1.297"
    "E8"  = "38.15%"
    "D9"  = "0.8714"
    "E9"  = "-0.87%"
    "E10" = "0.71%"
    "E11" = "-1.67%"
    "D12" = "0.07380"
    "E12" = "1.31%"
    "D13" = "0.02996"
    "E13" = "-5.25%"
    "D14" = "0.09050"
    "E14" = "-0.04%"
    "D15" = "0.001578"
    "E15" = "1.04%"
    "D16" = "0.0006329"
    "E16" = "0.83%"
    "D17" = "0.005911"
    "E17" = "-1.22%"
    "D18" = "3.454"
    "D19" = "2.285"
    "E19" = "-0.07%"
    "D21" = "0.1325"
    "E21" = "1.11%"
    "D22" = "3.904"
    "E22" = "1.52%"
    "D23" = "0.04360"
    "E23" = "1.00%"
    "D24" = "0.001175"
    "E24" = "0.00%"
    "D25" = "0.004262"
    "E25" = "-0.37%"
    "E26" = "0.09%"
    "D27" = "0.0001693"
    "E27" = "0.28%"
    "E40" = "0.05%"
    "D41" = "0.006685"
    "E41" = "4.30%"
    "E42" = "1.20%"
    "D43" = "0.002090"
    "E43" = "-0.65%"
    "D44" = "0.01217"
    "E44" = "-12.71%"
    "D45" = "0.00005309"
    "E45" = "3.06%"
    "D46" = "1.454"
    "E46" = "-38.04%"
    "D47" = "0.02000"
    "E47" = "-33.04%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
